$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: add "Co-Promotor" note in Found at/Journal column (M)
$ws.Range("M45").Value = "Co-Promotor"

# Row 46: O'Reilly - The Philosophy of Residuality Theory
$ws.Range("B46").Value = "None"
$ws.Range("C46").Value = "Antifragility"
$ws.Range("D46").Value = "The Philosophy of Residuality Theory"
$ws.Range("F46").Value = "Barry M O'Reilly"
$ws.Range("G46").Value = 2021
$ws.Range("H46").Value = "Journal Article"
$ws.Range("I46").Value = "M'OReilly2021"
$ws.Range("M46").Value = "ANT 2021"
$ws.Range("N46").Value = "10.1016/j.procs.2021.03.101"
$ws.Range("O46").Value = "https://doi.org/10.1016/j.procs.2021.03.101"
$ws.Range("P46").Value = "September, 2021"
$ws.Range("U46").Value = "Source received by Co-Promotor"

# Row 47: O'Reilly - An Introduction to Residuality Theory
$ws.Range("B47").Value = "None"
$ws.Range("C47").Value = "Antifragility"
$ws.Range("D47").Value = "An Introduction to Residuality Theory"
$ws.Range("E47").Value = "Software Design Heuristics for Complex Systems"
$ws.Range("F47").Value = "Barry M O'Reilly"
$ws.Range("G47").Value = 2020
$ws.Range("H47").Value = "Journal Article"
$ws.Range("I47").Value = "O’Reilly2020"
$ws.Range("M47").Value = "ANT 2020"
$ws.Range("N47").Value = "10.1016/j.procs.2020.03.120"
$ws.Range("O47").Value = "https://doi.org/10.1016/j.procs.2020.03.120"
$ws.Range("P47").Value = "September, 2021"
$ws.Range("U47").Value = "Source received by Co-Promotor"

# Row 48: Russo / Ciancarini - Towards Antifragile Software Architectures
$ws.Range("B48").Value = "None"
$ws.Range("C48").Value = "Antifragility"
$ws.Range("D48").Value = "Towards Antifragile Software Architectures"
$ws.Range("F48").Value = "Daniel Russo / Paolo Ciancarini"
$ws.Range("G48").Value = 2017
$ws.Range("H48").Value = "Journal Article"
$ws.Range("I48").Value = "Russo2017"
$ws.Range("M48").Value = "ANT2017"
$ws.Range("N48").Value = "10.1016/j.procs.2017.05.426"
$ws.Range("O48").Value = "https://doi.org/10.1016/j.procs.2017.05.426"
$ws.Range("P48").Value = "September, 2021"
$ws.Range("U48").Value = "Source received by Co-Promotor"

# Row 49: Russo / Ciancarini - A Proposal for an Antifragile Software Manifesto
$ws.Range("B49").Value = "None"
$ws.Range("C49").Value = "Antifragility"
$ws.Range("D49").Value = "A Proposal for an Antifragile Software Manifesto"
$ws.Range("F49").Value = "Daniel Russo / Paolo Ciancarini"
$ws.Range("G49").Value = 2016
$ws.Range("H49").Value = "Journal Article"
$ws.Range("I49").Value = "Russo2016"
$ws.Range("M49").Value = "ANT2016"
$ws.Range("N49").Value = "10.1016/j.procs.2016.04.196"
$ws.Range("O49").Value = "https://doi.org/10.1016/j.procs.2016.04.196"
$ws.Range("P49").Value = "September, 2021"
$ws.Range("U49").Value = "Source received by Co-Promotor"

# Row 50: Verhulst - Applying systems and safety engineering principles for antifragility
$ws.Range("B50").Value = "None"
$ws.Range("C50").Value = "Antifragility"
$ws.Range("D50").Value = "Applying systems and safety engineering principles for antifragility"
$ws.Range("F50").Value = "Eric Verhulst"
$ws.Range("G50").Value = 2014
$ws.Range("H50").Value = "Journal Article"
$ws.Range("I50").Value = "Verhulst2014"
$ws.Range("M50").Value = "ANT2014 Program: http://cs-conferences.acadiau.ca/ant-14/subPages/2014_Advanced_Program.pdf"
$ws.Range("N50").Value = "10.1016/j.procs.2014.05.500"
$ws.Range("O50").Value = "https://doi.org/10.1016/j.procs.2014.05.500"
$ws.Range("P50").Value = "September, 2021"
$ws.Range("U50").Value = "Source received by Co-Promotor"

# Row 51: Jones - Engineering Antifragile Systems
$ws.Range("B51").Value = "None"
$ws.Range("C51").Value = "Antifragility"
$ws.Range("D51").Value = "Engineering Antifragile Systems"
$ws.Range("E51").Value = "A Change In Design Philosophy"
$ws.Range("F51").Value = "Kennie H. Jones"
$ws.Range("G51").Value = 2014
$ws.Range("H51").Value = "Journal Article"
$ws.Range("I51").Value = "Jones2014"
$ws.Range("M51").Value = "ANT2014 Program: http://cs-conferences.acadiau.ca/ant-14/subPages/2014_Advanced_Program.pdf"
$ws.Range("N51").Value = "10.1016/j.procs.2014.05.504"
$ws.Range("O51").Value = "https://doi.org/10.1016/j.procs.2014.05.504"
$ws.Range("P51").Value = "September, 2021"
$ws.Range("U51").Value = "Source received by Co-Promotor"
